# Add a "Collection method" column to the pH water sampling log.
# The new column is inserted immediately before the existing
# "Sample Bottle Size (mL)" column (column N), and the previously
# trailing "Deployment Notes" column (which only ever repeated one
# boilerplate note) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- restructure columns -------------------------------------------------
# Insert a new blank column at N; everything from N onward shifts right.
$ws.Columns("N:N").Insert()

# The old last column (originally Z, "Deployment Notes") is now at AA.
# Remove it completely - its text is not kept anywhere else.
$ws.Columns("AA:AA").Delete()

# --- new column header ----------------------------------------------------
$ws.Range("N1").Value = "Collection method"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").WrapText = $true
$ws.Range("N1").Borders.LineStyle = 1
$ws.Columns("N:N").ColumnWidth = 10.0

# --- new column values ------------------------------------------------------
# rows 2-5 were sampled via the rosette; rows 6-13 via a single niskin
# (mirrors the existing "Sample notes" text, but using the normalized
# collection_method tokens)
$ws.Cells.Item(2, 14).Value = "rosette"
$ws.Cells.Item(3, 14).Value = "rosette"
$ws.Cells.Item(4, 14).Value = "rosette"
$ws.Cells.Item(5, 14).Value = "rosette"
$ws.Cells.Item(6, 14).Value = "single_niskin"
$ws.Cells.Item(7, 14).Value = "single_niskin"
$ws.Cells.Item(8, 14).Value = "single_niskin"
$ws.Cells.Item(9, 14).Value = "single_niskin"
$ws.Cells.Item(10, 14).Value = "single_niskin"
$ws.Cells.Item(11, 14).Value = "single_niskin"
$ws.Cells.Item(12, 14).Value = "single_niskin"
$ws.Cells.Item(13, 14).Value = "single_niskin"

# --- selection / view state -------------------------------------------------
$ws.Range("E24").Select()
